# Apply updated cryptocurrency price/volume figures per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'68.413.64"
$ws.Range("E2").Value = "  -6.81%  "

# Row 3
$ws.Range("D3").Value = "'3.741.66"
$ws.Range("E3").Value = "  -5.91%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'583.47"
$ws.Range("E5").Value = "  -4.79%  "

# Row 6
$ws.Range("D6").Value = "'177.38"
$ws.Range("E6").Value = "  +5.03%  "

# Row 7
$ws.Range("D7").Value = "'3.735.60"
$ws.Range("E7").Value = "  -5.85%  "

# Row 8
$ws.Range("D8").Value = "'0.635"
$ws.Range("E8").Value = "  -6.75%  "

# Row 9
$ws.Range("E9").Value = "  +0.27%  "

# Row 10
$ws.Range("D10").Value = "'0.721"
$ws.Range("E10").Value = "  -5.83%  "

# Row 11
$ws.Range("E11").Value = "  -10.15%  "

# Row 12
$ws.Range("D12").Value = "'54.11"
$ws.Range("E12").Value = "  -3.26%  "

# Row 13
$ws.Range("D13").Value = "'0.0000301"
$ws.Range("E13").Value = "  -10.94%  "

# Row 14
$ws.Range("D14").Value = "'10.81"
$ws.Range("E14").Value = "  -3.02%  "

# Row 15
$ws.Range("D15").Value = "'4.332.57"
$ws.Range("E15").Value = "  -6.15%  "

# Row 16
$ws.Range("D16").Value = "'3.768.56"
$ws.Range("E16").Value = "  -5.56%  "

# Row 17
$ws.Range("D17").Value = "'19.54"
$ws.Range("E17").Value = "  -4.42%  "

# Row 18
$ws.Range("E18").Value = "  -6.77%  "

# Row 19
$ws.Range("E19").Value = "  -6.74%  "

# Row 20
$ws.Range("E20").Value = "  -2.73%  "

# Row 21
$ws.Range("D21").Value = "'68.312.13"
$ws.Range("E21").Value = "  -6.87%  "

# Row 22
$ws.Range("D22").Value = "'413.26"
$ws.Range("E22").Value = "  -6.02%  "

# Row 23
$ws.Range("D23").Value = "'4.58"
$ws.Range("E23").Value = "  -5.74%  "

# Row 24
$ws.Range("D24").Value = "'89.25"
$ws.Range("E24").Value = "  -6.83%  "

# Row 25
$ws.Range("E25").Value = "  -7.71%  "

# Row 26
$ws.Range("D26").Value = "'12.99"
$ws.Range("E26").Value = "  -8.50%  "

# Row 27
$ws.Range("D27").Value = "'10.90"
$ws.Range("E27").Value = "  -1.38%  "

# Row 28
$ws.Range("D28").Value = "'3.88"
$ws.Range("E28").Value = "  -4.74%  "

# Row 29
$ws.Range("E29").Value = "  +0.46%  "

# Row 30
$ws.Range("D30").Value = "'9.64"
$ws.Range("E30").Value = "  -8.32%  "

# Row 31
$ws.Range("D31").Value = "'8.04"
$ws.Range("E31").Value = "  +3.41%  "

# Row 32
$ws.Range("D32").Value = "'33.27"
$ws.Range("E32").Value = "  -7.79%  "

# Row 33
$ws.Range("D33").Value = "'12.81"
$ws.Range("E33").Value = "  -7.46%  "

# Row 34
$ws.Range("E34").Value = "  -7.97%  "

# Row 35
$ws.Range("D35").Value = "'66.17"
$ws.Range("E35").Value = "  -6.30%  "

# Row 36
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").Value = "'616.18"
$ws.Range("E36").Value = "  -4.77%  "

# Row 37
$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "'44.15"
$ws.Range("E37").Value = "  -7.69%  "

# Row 38
$ws.Range("D38").Value = "'0.0" + [char]0x2083 + "0929"
$ws.Range("E38").Value = "  -12.96%  "

# Row 39
$ws.Range("D39").Value = "'0.403"
$ws.Range("E39").Value = "  -6.20%  "

# Row 40
$ws.Range("E40").Value = "  +0.21%  "

# Row 41
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.33%  "

# Row 42
$ws.Range("D42").Value = "'3.22"
$ws.Range("E42").Value = "  -1.38%  "

# Row 43
$ws.Range("E43").Value = "  -5.76%  "

# Row 44
$ws.Range("D44").Value = "'3.10"
$ws.Range("E44").Value = "  -8.56%  "

# Row 45
$ws.Range("E45").Value = "  -7.85%  "

# Row 46
$ws.Range("D46").Value = "'2.63"
$ws.Range("E46").Value = "  +2.63%  "

# Row 47
$ws.Range("D47").Value = "'9.47"
$ws.Range("E47").Value = "  -10.60%  "

# Row 48
$ws.Range("E48").Value = "  -8.16%  "

# Row 49
$ws.Range("D49").Value = "'2.73"
$ws.Range("E49").Value = "  -15.29%  "

# Row 50
$ws.Range("D50").Value = "'3.17"
$ws.Range("E50").Value = "  -7.35%  "

# Row 51
$ws.Range("D51").Value = "'2.735.73"
$ws.Range("E51").Value = "  -2.25%  "
